$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.018.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.67%  "
$ws.Range("D3").Value = "'1.870.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "'0.4329"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.86%  "
$ws.Range("E8").Value = "  -2.89%  "
$ws.Range("D9").Value = "'0.07399"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.76%  "
$ws.Range("D10").Value = "'0.9300"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.20%  "
$ws.Range("D11").Value = "'21.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.79%  "
$ws.Range("D12").Value = "'1.939.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "'6.730"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.53%  "
$ws.Range("D14").Value = "'5.424"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.76%  "
$ws.Range("D15").Value = "'0.06870"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.33%  "
$ws.Range("D16").Value = "'1.007"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "'80.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.84%  "
$ws.Range("D18").Value = "'0.000008981"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.80%  "
$ws.Range("D19").Value = "'1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("E20").Value = "  -6.01%  "
$ws.Range("D21").Value = "'28.004.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.75%  "
$ws.Range("D22").Value = "'5.116"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.22%  "
$ws.Range("D23").Value = "'10.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'2.176.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").Value = "'2.051"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("E26").Value = "  -2.11%  "
$ws.Range("D27").Value = "'18.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("D28").Value = "'5.481"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.37%  "
$ws.Range("D29").Value = "'113.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.09%  "
$ws.Range("D30").Value = "'1.685"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.10%  "
$ws.Range("D31").Value = "'0.08971"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("D32").Value = "'0.8041"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.49%  "
$ws.Range("D33").Value = "'4.754"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.00%  "
$ws.Range("D34").Value = "'1.172"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'2.954"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").Value = "'1.007"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").Value = "'0.05493"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.50%  "
$ws.Range("E38").Value = "  -3.53%  "
$ws.Range("D39").Value = "'0.01969"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("D40").Value = "'2.995"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.54%  "
$ws.Range("D41").Value = "'0.5231"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.13%  "
$ws.Range("D42").Value = "'7.025"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.39%  "
$ws.Range("D43").Value = "'0.1682"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.10%  "
$ws.Range("D44").Value = "'8.728"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.91%  "
$ws.Range("D45").Value = "'0.06708"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.98%  "
$ws.Range("D46").Value = "'0.4864"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.43%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.45%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'106.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.48%  "
$ws.Range("D49").Value = "'1.005"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "'1.668"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.66%  "
$ws.Range("D51").Value = "'1.864"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -15.00%  "
